$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new product row at row 21 (shifts the existing row 21 "سرنجات 3 سم"
# down to row 22, and row 22 "معجون اسنان فلورو بالكولا" down to row 23).
$ws.Rows("21:21").Insert()

# Copy formatting (styles only) from the row right below (the shifted-down old
# row 21, which already carries the correct per-column styles) into the new row.
$ws.Range("A22:N22").Copy()
$ws.Range("A21:N21").PasteSpecial(-4122)

# Restore the row height that row 21 had before the insert (Excel gives newly
# inserted rows the default height otherwise).
$ws.Rows("21:21").RowHeight = 24.75

# Recreate the merged-cell groups for the new row (Insert() does not carry
# merges over onto the freshly inserted blank row).
$ws.Range("B21:G21").Merge()
$ws.Range("H21:K21").Merge()
$ws.Range("L21:M21").Merge()

# Populate the new product row.
$ws.Range("A21").Value = 18
$ws.Range("B21").Value = "بيبي جوي رقم 4"
$ws.Range("H21").Value = "1:0"
$ws.Range("L21").Value = 320
$ws.Range("N21").Value = "1:0"

# Renumber the sequence column for the two rows that shifted down.
$ws.Range("A22").Value = 19
$ws.Range("A23").Value = 20

# Update the running total (was 1218.04, now increased by the new row's 320).
$ws.Range("K24").Value = 1538.04

# Minor row-height re-fit Excel applies to the footer row after the insert.
$ws.Rows("25:25").RowHeight = 16.5
